$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16: add a "control" header in column C ---
$ws.Range("C16").Value = "control"

# --- Row 18: add "List1" grouping label in column C ---
$ws.Range("C18").Value = "List1"

# --- Row 19: add "List2" grouping label in column C ---
$ws.Range("C19").Value = "List2"

# --- Row 20: add "List2" grouping label in column C ---
$ws.Range("C20").Value = "List2"

# --- Row 21 (new): repeat of the B4Ta tumor sample ---
$ws.Range("B21").Value = "s_B4Ta"

# --- Row 22 (new): repeat of the B4N normal sample ---
$ws.Range("A22").Value = "s_B4N"

# Remove the old row 24 (its content is being relocated to row 33 below)
$ws.Rows(24).Delete()

# --- Row 25 (new): [ControlPanel] section header ---
$ws.Range("A25").Value = "[ControlPanel]"
$ws.Range("A25").WrapText = $true

# --- Row 26 (new): List1 definition ---
$ws.Range("A26").Value = "List1"
$ws.Range("A26").WrapText = $true
$ws.Range("B26").Value = "s_B1N,s_B2Na, s_B2Nb, s_B3N, s_B4N"

# --- Row 27 (new): List2 definition ---
$ws.Range("A27").Value = "List2"
$ws.Range("A27").WrapText = $true
$ws.Range("B27").Value = "s_B3N, s_B4N"

# --- Row 33 (new): relocated comment/header row ---
$ws.Range("A33").Value = "#s_B1N  s_B1T  s_B2Na  s_B2Nb  s_B2T  s_B3N  s_B3Ta  s_B3Tb  s_B4N  s_B4Ta  s_B4Tb  s_B5N  s_B6T  s_B7T"

# Restore the selection to the cell the author left active
$ws.Range("A23").Select()
